# ChangeLoginPassword.xlsx -- "Limit Management Slider Update work"
#
# Adds a new "customer_password_query" column (R) with its SQL text, widens
# that column to fit the long query text, and moves the sheet's scroll/
# selection over to the new column.
#
# NOTE: the source diff also renamed the recorded Excel "last saved from"
# folder (xl/workbook.xml -> x15ac:absPath). That value is written by real
# Excel from the OS save path and is not exposed anywhere on the Excel COM
# object model (no settable property on Application/Workbook reaches it),
# so it cannot be reproduced from script here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column R: header + value ------------------------------------------
$ws.Range("R1").Value = "customer_password_query"
$ws.Range("R1").NumberFormat = "@"

$ws.Range("R2").Value = "Select count(*) from (Select MAX(AA.CREATED_ON) from DC_CUSTOMER_PASSWORD_HISTORY AA INNER JOIN DC_CUSTOMER_INFO BB ON AA.CUSTOMER_INFO_ID = BB.CUSTOMER_INFO_ID where AA.CUSTOMER_INFO_ID = (Select CUSTOMER_INFO_ID from dc_customer_info l where L.CUSTOMER_NAME = '{customer_name}') and AA.TRANSACTION_TYPE_ID = (Select LL.TRANSACTION_TYPE_ID from DC_TRANSACTION LL where LL.TRANSACTION_ID = '{TRANSACTION_ID}') and AA.PASSWORD = BB.LOGIN_PASSWORD and TRUNC(AA.CREATED_ON) < (SELECT TRUNC(SYSDATE) FROM DUAL) and TRUNC(AA.UPDATED_ON) < (SELECT TRUNC(SYSDATE) FROM DUAL) order by AA.UPDATED_ON desc) where rownum = 1"
$ws.Range("R2").Style = "Normal"

# Widen column R (bestFit-style) so the long query text is fully visible,
# splitting it out of the former R:S 9.140625-wide pair.
$ws.Columns.Item(18).ColumnWidth = 254.83333333333334

# --- Scroll / selection over to the new column ------------------------------
$ws.Activate()
$ws.Range("R11").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 18
$win.ScrollRow = 1
